$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.065182100734830328
$ws.Range("A2").Value = -0.0099999994890822563
$ws.Range("A3").Value = -0.0089999994932075111
$ws.Range("A4").Value = 0.28399311920450998
$ws.Range("A5").Value = -0.0059999995040236342
$ws.Range("A6").Value = -0.0059999994873400908
$ws.Range("A7").Value = -0.019999999410865499
$ws.Range("A8").Value = -0.019999999406985935
$ws.Range("A9").Value = -0.0059999994782229393
$ws.Range("A10").Value = 0.026510288155520811
$ws.Range("A11").Value = -0.0044999994821957046
$ws.Range("A12").Value = -0.0059999994729649231
$ws.Range("A13").Value = -0.0059999994656347866
$ws.Range("A14").Value = -0.011999999431347774
$ws.Range("A15").Value = -0.0059999994620341113
$ws.Range("A16").Value = -0.0059999994606236839
$ws.Range("A17").Value = -0.0059999994587425221
$ws.Range("A18").Value = -0.0089999994423592966
$ws.Range("A19").Value = -0.0089999994942346895
$ws.Range("A20").Value = -0.0089999994901948099
$ws.Range("A21").Value = -0.0089999994896521329
$ws.Range("A22").Value = -0.0089999994893625868
$ws.Range("A23").Value = -0.0089999994895384461
$ws.Range("A24").Value = -0.041999999304938918
$ws.Range("A25").Value = -0.041999999301569169
$ws.Range("A26").Value = -0.0059999994857271588
$ws.Range("A27").Value = -0.0059999994834747383
$ws.Range("A28").Value = -0.005999999472202866
$ws.Range("A29").Value = -0.022791407705152977
$ws.Range("A30").Value = -0.024256872543955588
$ws.Range("A31").Value = -0.014999999408278342
$ws.Range("A32").Value = -0.020999999375437284
$ws.Range("A33").Value = -0.0059999994557982106
